# RPA datasets push 2023-10-18
# A new IPO entry ("그린리소스") is inserted as row 3 on the
# "02_38커뮤니케이션(최근일자기준)" sheet, pushing the existing rows 3-11
# down to rows 4-12. The row that previously occupied the bottom of that
# block (the old row 12) is dropped, since it held stale/placeholder data
# that the refreshed row 3 now replaces.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Column order: A=회사명, B=수요예측일, C=희망공모가(원), D=확정공모가, E=공모금액(백만), F=주간사

$rows = @(
    @("와이바이오로직스", "2023.11.10~11.16", "9,000~11,000", "-", 13500, "유안타증권"),
    @("그린리소스",       "2023.11.03~11.09", "11,000~14,000", "-", 18040, "NH투자증권"),
    @("한선엔지니어링",   "2023.11.02~11.08", "5,200~6,000",   "-", 22100, "대신증권"),
    @("에코아이",         "2023.11.01~11.07", "28,500~34,700", "-", 59251, "KB증권"),
    @("동인기연(유가)",   "2023.11.01~11.07", "33,000~37,000", "-", 60654, "NH투자증권"),
    @("스톰테크",         "2023.10.31~11.06", "8,000~9,500",   "-", 26800, "하이투자증권"),
    @("블루엠텍",         "2023.10.31~11.06", "15,000~19,000", "-", 21000, "하나증권,키움증권"),
    @("에코프로머티리얼즈", "2023.10.30~11.03", "36,200~44,000", "-", 524031, "미래에셋증권,NH투자증권,하이투자증권"),
    @("캡스톤파트너스",   "2023.10.26~11.01", "3,200~3,600",   "-", 5107, "NH투자증권"),
    @("에이텀",           "2023.10.26~11.01", "23,000~30,000", "-", 14950, "하나증권"),
    @("한국스팩13호",     "2023.10.25~10.26", "2,000~2,000",   "-", 8000, "한국투자증권"),
    @("에이직랜드",       "2023.10.23~10.27", "19,100~21,400", "-", 50353, "삼성증권"),
    @("에스와이스틸텍",   "2023.10.23~10.27", "1,200~1,500",   "-", 8400, "KB증권"),
    @("컨텍",             "2023.10.20~10.26", "20,300~22,500", "-", 41818, "대신증권"),
    @("큐로셀",           "2023.10.20~10.26", "29,800~33,500", "-", 47680, "미래에셋증권,삼성증권"),
    @("메가터치",         "2023.10.20~10.26", "3,500~4,000",   "-", 18200, "NH투자증권"),
    @("비아이매트릭스",   "2023.10.19~10.25", "9,100~11,000",  "-", 10920, "IBK투자증권"),
    @("KB스팩27호",       "2023.10.19~10.20", "-",             "-", 25000, "KB증권"),
    @("유투바이오",       "2023.10.18~10.19", "3,300~3,900",   "-", 3724, "신한투자증권"),
    @("쏘닉스",           "2023.10.17~10.23", "5,000~7,000",   "-", 18000, "KB증권")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
